$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1214
$ws1.Range("F3").Value = 649
$ws1.Range("F5").Value = 5039
$ws1.Range("F7").Value = 9520
$ws1.Range("F8").Value = 245
$ws1.Range("F9").Value = 531
$ws1.Range("F11").Value = 678

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 20

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 649
$ws4.Range("F4").Value = 351
$ws4.Range("F7").Value = 5039
$ws4.Range("F10").Value = 9520
$ws4.Range("F11").Value = 245
$ws4.Range("F12").Value = 0
$ws4.Range("F13").Value = 88
$ws4.Range("F16").Value = 678
$ws4.Range("F18").Value = 0
